$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45181 -> 45182) for every data row (rows 2 through 375).
$ws.Range("C2:C375").Value = 45182
